# Add a new citation row for the Correlates of War Formal Alliances
# dataset (gibler2009ima) right after the existing "alliance" / ATOP row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ps_data_version")

# Insert a new row at position 6 (shifts existing rows 6-32 down to 7-33)
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the alliance citation data
$ws.Range("A6").Value = "alliance"
$ws.Range("B6").Value = "Correlates of War Formal Alliances"
$ws.Range("C6").Value = 4.1
$ws.Range("D6").Value = "gibler2009ima"

# Match the active cell selection recorded in the saved workbook
$ws.Range("M8").Select()
